$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bionomics")

# New "Tabled?" column in column A, header + "y" markers for the rows
# that represent directly-tabled (non-relational) fields.
$ws.Range("A2").Value = "Tabled?"

$taggedRows = @(4,6,7,8,9,10,12,13,14,15,16,17,18,19,20,21,23,24,25,26)
foreach ($r in $taggedRows) {
    $ws.Range("A$r").Value = "y"
}

# Move the active tab / selection from Species (last sheet) back to the
# first sheet, Bionomics, with the new selection at B22.
$ws.Activate()
$ws.Range("B22").Select()
